# Update countries & provincias Spain
# - Rusia overtakes Brasil (rows 5-6)
# - Estonia overtakes El Salvador (rows 90-91)
# - Eslovaquia overtakes Nueva Zelanda (rows 97-98)
# - Lituania row (94) case counts refreshed
# - Header timestamp bumped from 09:35 to 10:05

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 10:05"

# Row 5: Rusia (was Brasil) with refreshed figures
$ws.Range("A5").Value = "Rusia"
$ws.Range("B5").Value = 335882
$ws.Range("C5").Value = 9434
$ws.Range("D5").Value = 107936
$ws.Range("E5").Value = 224558
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 139
$ws.Range("H5").Value = 3388

# Row 6: Brasil (was Rusia), takes over the previous Brasil figures
$ws.Range("A6").Value = "Brasil"
$ws.Range("B6").Value = 332382
$ws.Range("C6").Value = 1492
$ws.Range("D6").Value = 135430
$ws.Range("E6").Value = 175836
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = 21116

# Row 90: Estonia (was El Salvador) with refreshed figures
$ws.Range("A90").Value = "Estonia"
$ws.Range("B90").Value = 1821
$ws.Range("C90").Value = 14
$ws.Range("D90").Value = 1526
$ws.Range("E90").Value = 231
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 64

# Row 91: El Salvador (was Estonia), takes over the previous El Salvador figures
$ws.Range("A91").Value = "El Salvador"
$ws.Range("B91").Value = 1819
$ws.Range("C91").Value = 94
$ws.Range("D91").Value = 570
$ws.Range("E91").Value = 1216
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 33

# Row 94: Lituania keeps its place, figures refreshed
$ws.Range("B94").Value = 1616
$ws.Range("C94").Value = 12
$ws.Range("D94").Value = 1135
$ws.Range("E94").Value = 418
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 63

# Row 97: Eslovaquia (was Nueva Zelanda) with refreshed figures
$ws.Range("A97").Value = "Eslovaquia"
$ws.Range("B97").Value = 1504
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1280
$ws.Range("E97").Value = 196
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 28

# Row 98: Nueva Zelanda (was Eslovaquia), takes over the previous Nueva Zelanda figures
$ws.Range("A98").Value = "Nueva Zelanda"
$ws.Range("B98").Value = 1504
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 1455
$ws.Range("E98").Value = 28
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 21
